$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 208239.92
$ws.Range("I19").Value = 387974.4
$ws.Range("J19").Value = 854
$ws.Range("K19").Value = 387974.4
$ws.Range("L19").Value = 854
$ws.Range("M19").Value = -387799.4
$ws.Range("N19").Value = -1204
$ws.Range("H41").Value = 154475.61
$ws.Range("I41").Value = 285878.84
$ws.Range("K41").Value = 285878.84
$ws.Range("M41").Value = -285438.84
$ws.Range("H113").Value = 3114.2727
$ws.Range("I113").Value = 3109.077
$ws.Range("J113").Value = 3121.7778
$ws.Range("K113").Value = 3109.077
$ws.Range("L113").Value = 3121.7778
$ws.Range("M113").Value = 144.9229999999998
$ws.Range("N113").Value = -9629.7778
$ws.Range("H116").Value = 3788.739
$ws.Range("I116").Value = 3427.1875
$ws.Range("J116").Value = 4615.143
$ws.Range("K116").Value = 3427.1875
$ws.Range("L116").Value = 4615.143
$ws.Range("M116").Value = 14.8125
$ws.Range("N116").Value = -11499.143
$ws.Range("H132").Value = 1314.6666
$ws.Range("I132").Value = 1252.9491
$ws.Range("J132").Value = 2225
$ws.Range("K132").Value = 3758.8473
$ws.Range("L132").Value = 6675
$ws.Range("M132").Value = -1228.8473
$ws.Range("N132").Value = -11735
$ws.Range("H137").Value = 850
$ws.Range("I137").Value = 700
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 2100
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = 450
$ws.Range("N137").Value = -8100

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2263.2222
$ws.Range("I2").Value = 1682.5333
$ws.Range("J2").Value = 5166.6665
$ws.Range("K2").Value = 1682.5333
$ws.Range("L2").Value = 5166.6665
$ws.Range("M2").Value = -1569.5333
$ws.Range("N2").Value = -5392.6665
$ws.Range("H6").Value = 10000000
$ws.Range("I6").Value = 10000000
$ws.Range("K6").Value = 10000000
$ws.Range("M6").Value = -9999827
$ws.Range("H32").Value = 12141
$ws.Range("I32").Value = 7228.7534
$ws.Range("J32").Value = 49801.555
$ws.Range("K32").Value = 7228.7534
$ws.Range("L32").Value = 49801.555
$ws.Range("M32").Value = -6941.7534
$ws.Range("N32").Value = -50375.555
$ws.Range("H39").Value = 11244.75
$ws.Range("I39").Value = 2489.5
$ws.Range("K39").Value = 2489.5
$ws.Range("M39").Value = -1969.5
$ws.Range("H45").Value = 1102.3334
$ws.Range("I45").Value = 1050
$ws.Range("J45").Value = 1207
$ws.Range("K45").Value = 1050
$ws.Range("L45").Value = 1207
$ws.Range("M45").Value = -673
$ws.Range("N45").Value = -1961
$ws.Range("H63").Value = 2253.7
$ws.Range("I63").Value = 2226.3333
$ws.Range("K63").Value = 2226.3333
$ws.Range("M63").Value = -1540.3333
$ws.Range("H66").Value = 2253.7
$ws.Range("I66").Value = 2226.3333
$ws.Range("K66").Value = 11131.6665
$ws.Range("M66").Value = -7699.666499999999
$ws.Range("H116").Value = 2263.2222
$ws.Range("I116").Value = 1682.5333
$ws.Range("J116").Value = 5166.6665
$ws.Range("K116").Value = 1682.5333
$ws.Range("L116").Value = 5166.6665
$ws.Range("M116").Value = 611.4666999999999
$ws.Range("N116").Value = -9754.666499999999
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H132").Value = 1464.2037
$ws.Range("I132").Value = 1277.5853
$ws.Range("J132").Value = 2052.7693
$ws.Range("K132").Value = 3832.7559
$ws.Range("L132").Value = 6158.3079
$ws.Range("M132").Value = -1302.7559
$ws.Range("N132").Value = -11218.3079

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2263.2222
$ws.Range("I3").Value = 1682.5333
$ws.Range("J3").Value = 5166.6665
$ws.Range("K3").Value = 1682.5333
$ws.Range("L3").Value = 5166.6665
$ws.Range("M3").Value = -1568.5333
$ws.Range("N3").Value = -5394.6665

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 6000149
$ws.Range("I12").Value = 6000149
$ws.Range("K12").Value = 6000149
$ws.Range("M12").Value = -5999979
$ws.Range("H16").Value = 1470.3334
$ws.Range("I16").Value = 1047.5714
$ws.Range("J16").Value = 2950
$ws.Range("K16").Value = 1047.5714
$ws.Range("L16").Value = 2950
$ws.Range("M16").Value = -760.5714
$ws.Range("N16").Value = -3524
$ws.Range("H70").Value = 14299.75
$ws.Range("J70").Value = 14299.75
$ws.Range("L70").Value = 14299.75
$ws.Range("N70").Value = -14929.75
$ws.Range("H73").Value = 14299.75
$ws.Range("J73").Value = 14299.75
$ws.Range("L73").Value = 14299.75
$ws.Range("N73").Value = -16483.75
$ws.Range("H86").Value = 21758
$ws.Range("I86").Value = 18596.666
$ws.Range("K86").Value = 18596.666
$ws.Range("M86").Value = -17473.666
$ws.Range("H89").Value = 21758
$ws.Range("I89").Value = 18596.666
$ws.Range("K89").Value = 92983.33
$ws.Range("M89").Value = -87367.33
$ws.Range("H113").Value = 1470.3334
$ws.Range("I113").Value = 1047.5714
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 1047.5714
$ws.Range("L113").Value = 2950
$ws.Range("M113").Value = 1122.4286
$ws.Range("N113").Value = -7290

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 150
$ws.Range("I40").Value = 150
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 600
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -531
$ws.Range("N40").ClearContents()
$ws.Range("H120").Value = 13618.182
$ws.Range("I120").Value = 11333.333
$ws.Range("J120").Value = 14475
$ws.Range("K120").Value = 33999.999
$ws.Range("L120").Value = 43425
$ws.Range("M120").Value = -29161.999
$ws.Range("N120").Value = -53101
$ws.Range("H131").Value = 17893776
$ws.Range("J131").Value = 1062.4783
$ws.Range("L131").Value = 3187.4349
$ws.Range("N131").Value = -13267.4349
$ws.Range("H132").Value = 1016.73334
$ws.Range("I132").Value = 876.9091
$ws.Range("J132").Value = 1401.25
$ws.Range("K132").Value = 7892.1819
$ws.Range("L132").Value = 12611.25
$ws.Range("M132").Value = -5362.1819
$ws.Range("N132").Value = -17671.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2750
$ws.Range("I43").Value = 2750
$ws.Range("K43").Value = 2750
$ws.Range("M43").Value = -2599
$ws.Range("H46").Value = 14666.667
$ws.Range("I46").Value = 6000
$ws.Range("J46").Value = 19000
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 19000
$ws.Range("M46").Value = -5844
$ws.Range("N46").Value = -19312
$ws.Range("H57").Value = 15527.5
$ws.Range("I57").Value = 10055
$ws.Range("J57").Value = 21000
$ws.Range("K57").Value = 10055
$ws.Range("L57").Value = 21000
$ws.Range("M57").Value = -9235
$ws.Range("N57").Value = -22640
$ws.Range("H63").Value = 11999.5
$ws.Range("J63").Value = 11999.5
$ws.Range("L63").Value = 11999.5
$ws.Range("N63").Value = -13371.5
$ws.Range("H66").Value = 11999.5
$ws.Range("J66").Value = 11999.5
$ws.Range("L66").Value = 35998.5
$ws.Range("N66").Value = -42862.5
$ws.Range("H132").Value = 2805.2683
$ws.Range("I132").Value = 2758.3333
$ws.Range("K132").Value = 8274.999899999999
$ws.Range("M132").Value = -5744.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 28000
$ws.Range("J123").Value = 28000
$ws.Range("L123").Value = 28000
$ws.Range("N123").Value = -37800
$ws.Range("H132").Value = 999.2069
$ws.Range("I132").Value = 1066.4
$ws.Range("J132").Value = 849.8889
$ws.Range("K132").Value = 3199.2
$ws.Range("L132").Value = 2549.6667
$ws.Range("M132").Value = -669.2000000000003
$ws.Range("N132").Value = -7609.6667
